$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "3631"

# Row 4
$ws.Range("A4").Value = "1073, 1105"
$ws.Range("C4").Value = "1105"
$ws.Range("D4").Value = "5131"

# Row 9
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "780, 780"
$ws.Range("D9").Value = "5677, 5887"

# Row 10
$ws.Range("A10").Value = "423, 1105, 1105"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "1105"
$ws.Range("D10").Value = "5331"

# Row 11
$ws.Range("A11").Value = "423, 1073, 1105"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "1105"
$ws.Range("D11").Value = "5433"

# Row 13
$ws.Range("A13").Value = "98, 130, 748, 1073"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "130, 130"
$ws.Range("D13").Value = "5582, 6488"

# Row 15
$ws.Range("A15").Value = "130, 423, 748, 1073"
$ws.Range("C15").Value = "130"
$ws.Range("D15").Value = "6561"

# Row 17
$ws.Range("A17").Value = "98, 780, 1073"
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = "1073, 98, 780"
$ws.Range("D17").Value = "6798, 6798, 6798"
